$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: Status Planned -> Confirmed
$ws.Range("D23").Value = "Confirmed"

# Row 24: Status Planned -> Confirmed
$ws.Range("D24").Value = "Confirmed"

# Row 39: Status Confirmed -> Delivered; Status update date 44309 -> 44316
$ws.Range("D39").Value = "Delivered"
$ws.Range("E39").Value = 44316

# Row 40: Status Planned -> Confirmed; Status update date 44302 -> 44316
$ws.Range("D40").Value = "Confirmed"
$ws.Range("E40").Value = 44316

# Row 58: Units 0 -> 1680
$ws.Range("B58").Value = 1680

# Row 60: Units 0 -> 4080
$ws.Range("B60").Value = 4080

# Row 61: Units 0 -> 1200
$ws.Range("B61").Value = 1200

# Row 62: Date 44320 -> 44316; Doses 50000 -> 180000; Status Assumption -> Delivered; Status update (blank) -> 44319
$ws.Range("A62").Value = 44316
$ws.Range("C62").Value = 180000
$ws.Range("D62").Value = "Delivered"
$ws.Range("E62").Value = 44319

# Row 63: Doses 50000 -> 33000; Status update (blank) -> 44312
$ws.Range("C63").Value = 33000
$ws.Range("E63").Value = 44312

# Row 72: Status Confirmed -> Delivered; Status update date 44309 -> 44316
$ws.Range("D72").Value = "Delivered"
$ws.Range("E72").Value = 44316
